# Fix mis-tagged property_category values:
#  - "建物" (building) sheet: property_category column (I), rows 2-12,
#    incorrectly says "land" -> should say "building"
#  - "汽車" (car) sheet: property_category column (H), rows 2-3,
#    incorrectly says "land" -> should say "car"

$wb = $excel.ActiveWorkbook

$wsBuilding = $wb.Worksheets.Item("建物")
for ($i = 2; $i -le 12; $i++) {
    $cell = $wsBuilding.Cells.Item($i, 9)
    $val = $cell.Value()
    if ($val -eq "land") {
        $cell.Value = "building"
    }
}

$wsCar = $wb.Worksheets.Item("汽車")
for ($i = 2; $i -le 3; $i++) {
    $cell = $wsCar.Cells.Item($i, 8)
    $val = $cell.Value()
    if ($val -eq "land") {
        $cell.Value = "car"
    }
}
